$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old array-formula results (swapTest / swapTest.multi) that
# used to live in E6:E10 and H4:H8, shrinking the used range back down.
$ws.Range("E6:E10").ClearContents()
$ws.Range("H4:H10").ClearContents()

# The curve-type/rate table that used to sit only at E3:G3 is now repeated
# down column E:G for rows 2,3,4,5,6 (multi-curve / multi-instrument test
# data). Row 3's "term" value changes from 2 to 3.
$ws.Range("E2").Value = "flatForward"
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 0.05

$ws.Range("F3").Value = 3

$ws.Range("E4").Value = "flatForward"
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 0.05

$ws.Range("E5").Value = "flatForward"
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 0.05

$ws.Range("E6").Value = "flatForward"
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 0.05

# Match the saved selection from the authored workbook.
[void]$ws.Range("F7").Select()
